$d = $word.ActiveDocument

# The paragraph currently repeats the "Question 1 - Correct" .. "Question 5 -
# Correct" block 5 times (each block separated by a blank line / double
# line-break). The edit keeps only the last occurrence of the block and
# removes the first four.

$block = "Question 1 - Correct" + [char]11 + "Question 2 - Correct" + [char]11 + "Question 3 - Correct" + [char]11 + "Question 4 - Correct" + [char]11 + "Question 5 - Correct" + [char]11 + [char]11

$blockLen = $block.Length
$removeCount = 4

$rng = $d.Range(0, $blockLen * $removeCount)
$rng.Delete()
